$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 22

# Column A holds a date-formatted text label (e.g. "2025-09-06") stored as
# plain text in the source file, not an actual date serial. Entering the
# string directly would trigger Excel's automatic date recognition, so we
# build it via a formula first and then convert the formula to its literal
# text result (copy / paste-special values), which preserves it as text.
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.Formula = '="2025-09-06"'
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)  # xlPasteValues

$ws.Cells.Item($row, 2).Value = 57.86999893188477
$ws.Cells.Item($row, 3).Value = 691.7000122070312
$ws.Cells.Item($row, 4).Value = 329.1499938964844
